$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.681.84'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').Value = '3.323.58'
$ws.Range('E3').Value = '  -4.62%  '
$ws.Range('E4').Value = '  +0.18%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '547.93'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.55%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '171.92'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -4.62%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.613'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -4.07%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '3.312.32'
$ws.Range('E9').Value = '  -4.52%  '
$ws.Range('E10').Value = '  -4.38%  '
$ws.Range('E11').Value = '  -1.17%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '53.22'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -2.32%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000265'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.10%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '8.85'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -4.64%  '
$ws.Range('D15').Value = '3.887.14'
$ws.Range('E15').Value = '  -3.81%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '18.18'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = '3.337.78'
$ws.Range('E17').Value = '  -4.22%  '
$ws.Range('E18').Value = '  -3.74%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '63.612.68'
$ws.Range('E19').Value = '  -3.04%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '11.62'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('E21').Value = '  -1.99%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '411.10'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.02'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '4.36'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.76%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '13.69'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +7.14%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '82.92'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.91%  '
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('E28').Value = '  -5.86%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '8.58'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -5.71%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '28.97'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.13%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.35'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -4.11%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '11.32'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -4.18%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '576.90'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -5.92%  '
$ws.Range('E34').Value = '  -4.28%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '57.67'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.147'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.36%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '34.93'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -7.47%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.41'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('E40').Value = '  -8.05%  '
$ws.Range('E41').Value = '  -5.04%  '
$ws.Range('D42').Value = '3.117.39'
$ws.Range('E42').Value = '  -7.33%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.40%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.77'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.57%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.23'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.98%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0399'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -4.35%  '
$ws.Range('E47').Value = '  -6.05%  '
$ws.Range('E48').Value = '  -4.32%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.127'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -4.08%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '132.60'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -3.80%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.99'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -5.77%  '
